$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat forced to
# text ("@") before assignment, otherwise Excel auto-converts them to numeric
# values (losing exact text representation such as trailing zeros).
$numericRiskCells = @(
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D13',
    'D14',
    'D15',
    'D17',
    'D18',
    'D19',
    'D20',
    'D22',
    'D23',
    'D27',
    'D28',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)

foreach ($cellRef in $numericRiskCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply new values
$ws.Range('D2').Value = '24.381.69'
$ws.Range('E2').Value = '  +9.80%  '
$ws.Range('D3').Value = '1.674.92'
$ws.Range('E3').Value = '  +5.48%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '0.9985'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').Value = '305.20'
$ws.Range('E6').Value = '  +2.37%  '
$ws.Range('D7').Value = '0.3683'
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('D8').Value = '0.3425'
$ws.Range('E8').Value = '  +2.71%  '
$ws.Range('D9').Value = '47.59'
$ws.Range('E9').Value = '  +15.66%  '
$ws.Range('D10').Value = '1.157'
$ws.Range('E10').Value = '  +4.01%  '
$ws.Range('D11').Value = '0.07204'
$ws.Range('E11').Value = '  +3.90%  '
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').Value = '6.128'
$ws.Range('E13').Value = '  +5.60%  '
$ws.Range('D14').Value = '20.09'
$ws.Range('E14').Value = '  +3.98%  '
$ws.Range('D15').Value = '6.719'
$ws.Range('E15').Value = '  +3.15%  '
$ws.Range('D16').Value = '1.674.76'
$ws.Range('E16').Value = '  +5.81%  '
$ws.Range('D17').Value = '0.00001099'
$ws.Range('E17').Value = '  +3.88%  '
$ws.Range('D18').Value = '0.9984'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '0.06645'
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').Value = '80.34'
$ws.Range('E20').Value = '  +5.88%  '
$ws.Range('E21').Value = '  +4.41%  '
$ws.Range('D22').Value = '6.093'
$ws.Range('E22').Value = '  +3.45%  '
$ws.Range('D23').Value = '12.15'
$ws.Range('E23').Value = '  +4.87%  '
$ws.Range('D24').Value = '24.323.34'
$ws.Range('E24').Value = '  +9.60%  '
$ws.Range('E26').Value = '  +6.58%  '
$ws.Range('D27').Value = '152.64'
$ws.Range('E27').Value = '  +3.59%  '
$ws.Range('D28').Value = '19.38'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').Value = '1.861.04'
$ws.Range('E29').Value = '  +6.23%  '
$ws.Range('D30').Value = '127.32'
$ws.Range('E30').Value = '  +5.36%  '
$ws.Range('D31').Value = '6.271'
$ws.Range('E31').Value = '  +7.44%  '
$ws.Range('D32').Value = '4.021'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').Value = '0.9686'
$ws.Range('E33').Value = '  +6.13%  '
$ws.Range('D34').Value = '0.08448'
$ws.Range('E34').Value = '  +4.01%  '
$ws.Range('D35').Value = '1.675'
$ws.Range('E35').Value = '  +3.17%  '
$ws.Range('D36').Value = '12.32'
$ws.Range('E36').Value = '  +5.80%  '
$ws.Range('D37').Value = '0.06377'
$ws.Range('E37').Value = '  +7.06%  '
$ws.Range('D38').Value = '5.290'
$ws.Range('E38').Value = '  +3.96%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '8.688'
$ws.Range('E39').Value = '  +4.42%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.02308'
$ws.Range('E40').Value = '  +5.75%  '
$ws.Range('D41').Value = '1.233'
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  +5.52%  '
$ws.Range('D43').Value = '0.6062'
$ws.Range('E43').Value = '  +5.15%  '
$ws.Range('D44').Value = '0.9983'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').Value = '3.746'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('D46').Value = '12.90'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').Value = '0.5855'
$ws.Range('E47').Value = '  +5.70%  '
$ws.Range('D48').Value = '125.38'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('D49').Value = '2.013'
$ws.Range('E49').Value = '  +4.34%  '
$ws.Range('D50').Value = '0.07137'
$ws.Range('E50').Value = '  +6.31%  '
$ws.Range('D51').Value = '75.67'
$ws.Range('E51').Value = '  +4.87%  '

# Restore default (Normal) cell style so no stray number-format styling remains
foreach ($cellRef in $numericRiskCells) {
    $ws.Range($cellRef).Style = "Normal"
}
